# Update the "修改时间" (modified time) column on each portfolio sheet to
# reflect the latest sync timestamp produced by the Web UI update.
$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-09-21 22:46:47"

# Sheet 1: 大智投资组合 -- timestamps live in column E, data rows 2-9
$ws1 = $wb.Worksheets.Item("大智投资组合")
$ws1.Range("E2:E9").Value = $newTimestamp

# Sheet 2: 大成投资组合 -- timestamps live in column E, data rows 2-11
$ws2 = $wb.Worksheets.Item("大成投资组合")
$ws2.Range("E2:E11").Value = $newTimestamp

# Sheet 3: 我的投资组合 -- timestamps live in column G, data rows 2-13
$ws3 = $wb.Worksheets.Item("我的投资组合")
$ws3.Range("G2:G13").Value = $newTimestamp
